$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Append two new rows to the end of the table (Table2), which currently
# spans A1:H153, growing it to A1:H155.
$lo.ListRows.Add() | Out-Null
$lo.ListRows.Add() | Out-Null

# --- Row 154: BoardMembership / Read / Refresh ---
$ws.Range("A154").Value = "BoardMembership"
$ws.Range("B154").Value = "Read"
$ws.Range("C154").Value = "Refresh"
$ws.Range("D154").Formula = "=Table2[[#This Row],[Entity]]&""_""&Table2[[#This Row],[R/W]]&""_""&Table2[[#This Row],[Requested Action]]"
$ws.Range("E154").Formula = "=Table2[[#This Row],[Enum Value]]&"",""" 
$ws.Range("F154").Value = """boards"",""_boardId"",""memberships"",""_id"""
$ws.Range("G154").Value = "Get"
$ws.Range("H154").Formula = "=""{EntityRequestType.""&Table2[[#This Row],[Enum Value]]&"", () => new Endpoint(RestMethod.""&Table2[[#This Row],[Method]]&"", new[]{""&Table2[[#This Row],[Endpoint]]&""})},"""

# --- Row 155: Attachment / Read / Refresh ---
$ws.Range("A155").Value = "Attachment"
$ws.Range("B155").Value = "Read"
$ws.Range("C155").Value = "Refresh"
$ws.Range("D155").Formula = "=Table2[[#This Row],[Entity]]&""_""&Table2[[#This Row],[R/W]]&""_""&Table2[[#This Row],[Requested Action]]"
$ws.Range("E155").Formula = "=Table2[[#This Row],[Enum Value]]&"",""" 
$ws.Range("F155").Value = """cards"",""_cardId"",""attachments"",""_id"""
$ws.Range("G155").Value = "Get"
$ws.Range("H155").Formula = "=""{EntityRequestType.""&Table2[[#This Row],[Enum Value]]&"", () => new Endpoint(RestMethod.""&Table2[[#This Row],[Method]]&"", new[]{""&Table2[[#This Row],[Endpoint]]&""})},"""

# Update the view: scroll the frozen (bottom) pane so row 125 is near the
# top, and select the last cell of the newly added data, matching where
# the author ended up after typing in the new rows.
$win = $excel.ActiveWindow
$win.ScrollRow = 125
$ws.Range("H155").Select()
